$wb = $excel.ActiveWorkbook

# 1. Remove the EU:US adjustment ratio from the ICtPSFfL incremental-cost
#    formulas (row 7, columns B:AK) before we remove the source cells the
#    ratio depended on (About!B27/About!B28).
$ws3 = $wb.Worksheets.Item("ICtPSFfL")
$cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD","AE","AF","AG","AH","AI","AJ","AK")
foreach ($col in $cols) {
    $ws3.Range($col + "7").Formula = "=MAX(Calcs!" + $col + "35,0)"
}

# 2. Remove the now-unused EU/US adjustment note and ratio rows from the
#    About sheet (rows 26:28).
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("A26:A28").EntireRow.Delete()
